$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 187; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = "BEST2_C_$i"
}

$ws.Range("B188").Select()
